$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added new test data: bump quantities for a few existing tasks
$ws.Range("C12").Value = 2
$ws.Range("C14").Value = 2
$ws.Range("C16").Value = 2

# Highlight the "Order Entry" row with a yellow fill
$ws.Range("A18").Interior.Color = 65535

# Column A grew a bit wider to fit the new data
$ws.Columns("A").ColumnWidth = 38.333333333333336

# Move the active selection
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
